# Bug fix for task2: the "instructionType" worksheet previously held a small
# (and stale/incorrect) 4-row summary of instruction-category counts. The fix
# replaces it with a full per-instruction mapping: column A = instruction
# mnemonic, column B = its parameter-count category - sourced straight from
# the already-correct "Summary" sheet (columns A and C), which enumerates all
# 164 instructions.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item("Summary")
$target  = $wb.Worksheets.Item("instructionType")

# Find how many data rows exist on the Summary sheet.
$lastRow = $summary.Cells.Item($summary.Rows.Count, 1).End(-4162).Row

# Clear the old 4-row aggregate table before laying down the new data.
$target.Cells.ClearContents()

for ($i = 1; $i -le $lastRow; $i++) {
    $name = $summary.Cells.Item($i, 1).Value2
    $cat  = $summary.Cells.Item($i, 3).Value2

    $target.Cells.Item($i, 1).Value2 = $name
    $target.Cells.Item($i, 2).Value2 = $cat
}
